# "added more wines to weinbar"
#
# Row 113 already held the "Marani Tvishi" wine but with several columns
# still empty - fill in the grape, region, vineyard, alcohol and extend the
# tasting-note description.
# Row 114 is a brand new wine ("Regent" from Michel Schneider / Edeka).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 113: flesh out the existing "Marani Tvishi" entry ---
$ws.Range("B113").Value = "Tsolikauri"
$ws.Range("E113").Value = "Racha-Lechkhumi"
$ws.Range("F113").Value = "Marani"
$ws.Range("G113").Value = "suess, herb, sueffig, banane"
$ws.Range("I113").Value = 11

# --- Row 114: add a new wine entry ---
$ws.Range("A114").Value = "Regent"
$ws.Range("B114").Value = "Regent"
$ws.Range("C114").Value = "red"
$ws.Range("D114").Value = "Germany"
$ws.Range("E114").Value = "Rheinhessen"
$ws.Range("F114").Value = "Michel Schneider"
$ws.Range("G114").Value = "nelken, schwer, zimt, flach, bitter"
$ws.Range("H114").Value = "Edeka"
$ws.Range("I114").Value = 12
$ws.Range("J114").Value = 1.99
$ws.Range("K114").Value = 2019
$ws.Range("L114").Value = "yes"
$ws.Range("M114").Value = "no"

# Column B ("grape") now has a long new value - let Excel auto-fit it.
$ws.Columns.Item(2).AutoFit() | Out-Null

# Leave the selection on the newly added cell, like the author did.
$ws.Range("F114").Select() | Out-Null
